# iron_native Word COM-interop script
# Applies the "latest update to how to" edit to howToVsCodeGit.docx:
#   1. Moves the "_GoBack" bookmark from the very end of the document to the
#      "Python console installed appropriately" bullet, and un-bullets that
#      paragraph (drops pStyle=ListParagraph / numPr).
#   2. Merges the "Yes" / " you have managed to open " runs (and drops the
#      surrounding proofErr gramStart/gramEnd markers) into a single run.
#   3. Merges the "Download the " / "64 bit" / " version of Git Windows" runs
#      (and drops the surrounding proofErr gramStart/gramEnd markers) into a
#      single run.
#   4. Removes the "_GoBack" bookmark that used to sit after
#      "You NOW have your own branch on git".

$d = $word.ActiveDocument

function Replace-ParagraphXml([object]$doc, [object]$paraRange, [string]$innerPXml) {
    # Replaces the contents of a single paragraph (identified by $paraRange,
    # typically obtained from Paragraphs(n).Range) with the supplied <w:p>...</w:p>
    # markup, preserving whatever lies before/after it in the document.
    #
    # Special-cased for the situation where $paraRange reaches all the way to
    # the end of the document's Content: Word can never truly delete the very
    # last paragraph mark of the body, so InsertXML on a range that includes
    # it actually inserts a new paragraph *before* that mark instead of
    # replacing in place. We detect that, and afterwards delete the extra
    # paragraph mark that gets left behind so the paragraph count returns to
    # normal.
    $isAtVeryEnd = ($paraRange.End -eq $doc.Content.End)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
           '<w:body>' + $innerPXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $beforeCount = $doc.Paragraphs.Count
    $paraRange.InsertXML($xml)

    if ($isAtVeryEnd -and $doc.Paragraphs.Count -gt $beforeCount) {
        # Collapse the spurious trailing empty paragraph back into the one we
        # just inserted by deleting the paragraph mark that separates them.
        $mergeEnd = $doc.Paragraphs($doc.Paragraphs.Count - 1).Range.End
        $doc.Range($mergeEnd - 1, $mergeEnd).Delete()
    }
}

# ---------------------------------------------------------------------------
# 1) "Python console installed appropriately": drop list formatting, add the
#    _GoBack bookmark right after the (now style-less) pPr.
# ---------------------------------------------------------------------------
$find = $d.Content
$null = $find.Find.Execute("Python console installed appropriately", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $find.Paragraphs(1).Range
$p1 = '<w:p w14:paraId="483D3414" w14:textId="63E2CB8D" w:rsidR="00E42783" w:rsidRPr="00E42783" w:rsidRDefault="00E42783" w:rsidP="00E42783">' + `
        '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
        '<w:r w:rsidRPr="00E42783"><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Python console installed appropriately</w:t></w:r>' + `
      '</w:p>'
Replace-ParagraphXml $d $target $p1

# ---------------------------------------------------------------------------
# 2) "Yes" + " you have managed to open " -> single run, proofErr removed.
# ---------------------------------------------------------------------------
$find = $d.Content
$null = $find.Find.Execute("Yes you have managed to open VSCode", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $find.Paragraphs(1).Range
$p2 = '<w:p w14:paraId="0780B63F" w14:textId="39759960" w:rsidR="00E42783" w:rsidRDefault="00E42783" w:rsidP="00E42783">' + `
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Yes you have managed to open </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>VSCode</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and the </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>machineLearningProject</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
      '</w:p>'
Replace-ParagraphXml $d $target $p2

# ---------------------------------------------------------------------------
# 3) "Download the " + "64 bit" + " version of Git Windows" -> single run,
#    proofErr removed.
# ---------------------------------------------------------------------------
$find = $d.Content
$null = $find.Find.Execute("Download the 64 bit version of Git Windows", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $find.Paragraphs(1).Range
$p3 = '<w:p w14:paraId="2FCB8EDB" w14:textId="15054FC2" w:rsidR="00E42783" w:rsidRDefault="00E42783" w:rsidP="00E42783">' + `
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Download the 64 bit version of Git Windows</w:t></w:r>' + `
      '</w:p>'
Replace-ParagraphXml $d $target $p3

# ---------------------------------------------------------------------------
# 4) Remove the _GoBack bookmark that used to trail
#    "You NOW have your own branch on git".
# ---------------------------------------------------------------------------
$find = $d.Content
$null = $find.Find.Execute("You NOW have your own branch on git", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $find.Paragraphs(1).Range
$p4 = '<w:p w14:paraId="0DD9A5B4" w14:textId="5022B14A" w:rsidR="00D76F78" w:rsidRPr="00E42783" w:rsidRDefault="00D76F78" w:rsidP="00E42783">' + `
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>You NOW have your own branch on git</w:t></w:r>' + `
      '</w:p>'
Replace-ParagraphXml $d $target $p4
